# Update the four "Dates de la campanya ..." paragraphs so the campaign
# dates sentence gets the new lead-in text, per the commit diff.
#
# We avoid Word's Find/Replace for the replacement text because this
# runtime's AutoFormat/AutoCorrect smart-quote substitution kicks in on
# Find.Execute's replacement string and turns the straight apostrophe in
# "d'Hèrcules" into a curly one, which does not match the target OOXML.
# Writing directly to the matching paragraph Range keeps the apostrophe
# (and the middot in "Constel·lació") intact.

$d = $word.ActiveDocument

$middot = [char]0x00B7

$old = "Dates de la campanya Constel" + $middot + "lació d'Hèrcules 2022: 13-22 de juny, 12-21 de juliol, 10-19 d'agost"
$new = "Dates de la campanya 2022 en què usem la constel" + $middot + "lació, Constel" + $middot + "lació d'Hèrcules 13-22 de juny, 12-21 de juliol, 10-19 d'agost"

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $t = $r.Text
    # Strip trailing paragraph mark (CR), cell mark (BEL) and the
    # section-break marker (FF) that Range.Text can append.
    $trimmed = $t.TrimEnd([char]13, [char]7, [char]12)
    if ($trimmed.Equals($old)) {
        $r.Text = $new
    }
}
